$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (reverting the "fucked it all up" commit)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2

$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 3

# Update the selected cell to match the saved selection state
$ws.Range("D9").Select()
